$d = $word.ActiveDocument

function Replace-NextOne([string]$findText, [string]$replaceText) {
    $r = $d.Content
    $ok = $r.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 0, $false, "", 0)
    if (-not $ok) {
        Write-Host "NOT FOUND: $findText"
        return
    }
    $r.Text = $replaceText
}

Replace-NextOne '2025-11-28' '2026-01-05'
Replace-NextOne '01:26 UTC' '01:19 UTC'
Replace-NextOne 'transactionId: ONTE176431099' 'transactionId: ONTE176758858'
Replace-NextOne '"ONTE176431099"' '"ONTE176758858"'
Replace-NextOne '"eyJhbGciOiJIUzI1NiIsInR5cCI6IkpXVCJ9.eyJtZXJjaGFudENvZGUiOiI0MDAxODM0IiwiZmFjaWxpdGF0b3JDb2RlIjoiMCIsInRyYW5zYWN0aW9uSWQiOiJPTlRFMTc2NDMxMDk5IiwiT3JkZXJOdW1iZXIiOiJPTlRFMTc2NDMxMDk5IiwiQW1vdW50IjoiMS4wMCIsIlRva2VuSWQiOiIyODRkN2MyMS0xZWFjLTQ3YjUtOWU1Yi04ZGM2N2QxMDkwMWYiLCJuYmYiOjE3NjQzMTEyMTYsImV4cCI6MTc2NDMxMjExNiwiaWF0IjoxNzY0MzExMjE2fQ.4EFPPYsehQc7dW_nbJmwDBppaZs8mjoAkcDXpp0Cl04"' '"eyJhbGciOiJIUzI1NiIsInR5cCI6IkpXVCJ9.eyJtZXJjaGFudENvZGUiOiI0MDAxODM0IiwiZmFjaWxpdGF0b3JDb2RlIjoiMCIsInRyYW5zYWN0aW9uSWQiOiIxNzY3NTkzOTQ4OTMxIiwiT3JkZXJOdW1iZXIiOiJPTlRFMTc2NzU4ODU4IiwiQW1vdW50IjoiMS4wMCIsIlRva2VuSWQiOiI3MTMyMmY3ZS0wNWJmLTRmMjUtYTEyNy0xOWM5Y2M2ZjNhNDIiLCJuYmYiOjE3Njc1OTM5NDksImV4cCI6MTc2NzU5NDg0OSwiaWF0IjoxNzY3NTkzOTQ5fQ.QgxElaAXjETc4-_hYAbj0xcYgYDLVJmNVLE-HsllGiM"'
Replace-NextOne 'transactionId: ONTE176431099' 'transactionId: ONTE176758858'
Replace-NextOne '"ONTE176431099"' '"ONTE176758858"'
Replace-NextOne '"S46451"' '"S16082"'
Replace-NextOne '"1546802"' '"1548093"'
Replace-NextOne '"2025-11-28 01:23:21.000"' '"2026-01-04 23:49:52.000"'
Replace-NextOne '"2025-11-28 01:23:21.000"' '"2026-01-04 23:49:52.000"'
Replace-NextOne '"2025-11-28 01:26:57.010"' '"2026-01-05 01:19:09.501"'
Replace-NextOne '"2025-11-28 01:26:57.208"' '"2026-01-05 01:19:09.824"'
Replace-NextOne '"198"' '"323"'
Replace-NextOne '"ONTE176431099"' '"ONTE176758858"'
Replace-NextOne '"card"' '"CARD"'
Replace-NextOne '"S46451"' '"S16082"'
Replace-NextOne '"1546802"' '"1548093"'
Replace-NextOne '"1243"' '"1271"'
Replace-NextOne '1243' '1271'
Replace-NextOne '"S46451"' '"S16082"'
Replace-NextOne '"0973473"' '"0861755"'
Replace-NextOne '"ONTE176431099"' '"ONTE176758858"'
Replace-NextOne '"20251128"' '"20260104"'
Replace-NextOne '"012321"' '"234952"'
Replace-NextOne '"AE"' '"MC"'
Replace-NextOne '"377753*****0152"' '"520474******1127"'
Replace-NextOne '"1546802"' '"1548093"'
Replace-NextOne '"pZmp6ZQc5ksCIrQvypJEoRAceGMDtH/YeWv9SrfR4ms="' '"pUUhF+5KrGfioTXa6YiJ7HIGn1zlQBXwM/iLrIYaFfw="'
Replace-NextOne 'S46451' 'S16082'
Replace-NextOne '0973473' '0861755'
Replace-NextOne 'ONTE176431099' 'ONTE176758858'
Replace-NextOne '20251128' '20260104'
Replace-NextOne '012321' '234952'
Replace-NextOne 'AE' 'MC'
Replace-NextOne '377753*****0152' '520474******1127'
Replace-NextOne '1546802' '1548093'
Replace-NextOne 'CARD (AE 377753*****0152)' 'CARD (MC 520474******1127)'
Replace-NextOne 'S46451' 'S16082'
Replace-NextOne '1546802' '1548093'
Replace-NextOne '1243' '1271'
Replace-NextOne 'en ambas operaciones (ONTE176431099)' 'en ambas operaciones (ONTE176758858)'
Replace-NextOne '(S46451)' '(S16082)'
Replace-NextOne '(1546802)' '(1548093)'
